$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-04 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-05 Saturday", 2) | Out-Null
$d.Content.Find.Execute("78×43=3354", $true, $false, $false, $false, $false, $true, 1, $false, "67×61=4087", 2) | Out-Null
$d.Content.Find.Execute("93×63=5859", $true, $false, $false, $false, $false, $true, 1, $false, "65×98=6370", 2) | Out-Null
$d.Content.Find.Execute("59×17=1003", $true, $false, $false, $false, $false, $true, 1, $false, "18×13=234", 2) | Out-Null
$d.Content.Find.Execute("27×41=1107", $true, $false, $false, $false, $false, $true, 1, $false, "53×24=1272", 2) | Out-Null
$d.Content.Find.Execute("99×89=8811", $true, $false, $false, $false, $false, $true, 1, $false, "22×91=2002", 2) | Out-Null
$d.Content.Find.Execute("85×85=7225", $true, $false, $false, $false, $false, $true, 1, $false, "94×14=1316", 2) | Out-Null
$d.Content.Find.Execute("88×47=4136", $true, $false, $false, $false, $false, $true, 1, $false, "88×87=7656", 2) | Out-Null
$d.Content.Find.Execute("53×32=1696", $true, $false, $false, $false, $false, $true, 1, $false, "36×27=972", 2) | Out-Null
$d.Content.Find.Execute("71×45=3195", $true, $false, $false, $false, $false, $true, 1, $false, "96×57=5472", 2) | Out-Null
$d.Content.Find.Execute("71×75=5325", $true, $false, $false, $false, $false, $true, 1, $false, "84×11=924", 2) | Out-Null
$d.Content.Find.Execute("59×90=5310", $true, $false, $false, $false, $false, $true, 1, $false, "67×22=1474", 2) | Out-Null
$d.Content.Find.Execute("14×96=1344", $true, $false, $false, $false, $false, $true, 1, $false, "21×33=693", 2) | Out-Null
$d.Content.Find.Execute("64×31=1984", $true, $false, $false, $false, $false, $true, 1, $false, "75×49=3675", 2) | Out-Null
$d.Content.Find.Execute("31×54=1674", $true, $false, $false, $false, $false, $true, 1, $false, "62×31=1922", 2) | Out-Null
$d.Content.Find.Execute("80×42=3360", $true, $false, $false, $false, $false, $true, 1, $false, "83×53=4399", 2) | Out-Null
$d.Content.Find.Execute("90×29=2610", $true, $false, $false, $false, $false, $true, 1, $false, "60×38=2280", 2) | Out-Null
$d.Content.Find.Execute("61×20=1220", $true, $false, $false, $false, $false, $true, 1, $false, "47×28=1316", 2) | Out-Null
$d.Content.Find.Execute("55×50=2750", $true, $false, $false, $false, $false, $true, 1, $false, "56×63=3528", 2) | Out-Null
$d.Content.Find.Execute("29×14=406", $true, $false, $false, $false, $false, $true, 1, $false, "20×52=1040", 2) | Out-Null
$d.Content.Find.Execute("14×13=182", $true, $false, $false, $false, $false, $true, 1, $false, "19×12=228", 2) | Out-Null
$d.Content.Find.Execute("50×65=3250", $true, $false, $false, $false, $false, $true, 1, $false, "41×96=3936", 2) | Out-Null
$d.Content.Find.Execute("71×66=4686", $true, $false, $false, $false, $false, $true, 1, $false, "64×52=3328", 2) | Out-Null
$d.Content.Find.Execute("58×28=1624", $true, $false, $false, $false, $false, $true, 1, $false, "97×33=3201", 2) | Out-Null
$d.Content.Find.Execute("57×19=1083", $true, $false, $false, $false, $false, $true, 1, $false, "62×17=1054", 2) | Out-Null
$d.Content.Find.Execute("53×71=3763", $true, $false, $false, $false, $false, $true, 1, $false, "93×46=4278", 2) | Out-Null
